# Weekly update: insert a new day's price record for Pomelo (Vega Modelo de
# Temuco) at the top of the data block (row 81), pushing all existing rows
# (old 81-132) down by one (new 82-133).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 81; rows below (old 81..132) shift to 82..133.
$ws.Rows.Item(81).Insert()

# Populate the new row 81 with this week's record. Most fields repeat the
# values that were already present for this market/product/variety/quality
# combination; only the date, volume, max/avg price and $/Kg columns change.
$ws.Range("A81").Value = 10
$ws.Range("B81").Value = "Vega Modelo de Temuco"
$ws.Range("C81").Value = "La Araucanía"
$ws.Range("D81").Value = 44452
$ws.Range("E81").Value = 9
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100102
$ws.Range("H81").Value = "Cítricos"
$ws.Range("I81").Value = 100102006
$ws.Range("J81").Value = "Pomelo"
$ws.Range("K81").Value = "Start Ruby"
$ws.Range("L81").Value = "Primera"
$ws.Range("M81").Value = 174
$ws.Range("N81").Value = 12000
$ws.Range("O81").Value = 13000
$ws.Range("P81").Value = 12431
$ws.Range("Q81").Value = "$/bandeja 15 kilos granel"
$ws.Range("R81").Value = "Región de O'Higgins"
$ws.Range("S81").Value = 829
$ws.Range("T81").Value = 15
